$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value in A2
$ws.Range("A2").Value = "6001-1122-STAFF-MJ0AJPEB"

# Set column A width to match the new content width (as if autofit after the edit)
$ws.Columns.Item(1).ColumnWidth = 33.3

# Update the active selection to B2, as seen in the target file
$ws.Range("B2").Select() | Out-Null
